$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old sector labels (GICS-style) to new sector labels (Morningstar-style)
$sectorMap = @{
    "Consumer Staples"       = "Consumer Defensive"
    "Consumer Discretionary" = "Consumer Cyclical"
    "Information Technology" = "Technology"
    "Other"                  = "Financial Services"
    "Materials"              = "Basic Materials"
}

# Find last used row in column D (sector column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Text
    if ($sectorMap.ContainsKey($current)) {
        $cell.Value = $sectorMap[$current]
    }
}
